$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to stay as Text (they hold numeric-looking strings
# like "574.37" that Excel would otherwise auto-convert to numbers);
# the source data models them as plain text, matching the original file.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.758.26"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.463.12"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D5").Value = "574.37"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "147.65"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "2.464.18"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "29.15"
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D16").Value = "2.912.63"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "62.699.67"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "2.469.71"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "7.94"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "10.97"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").Value = "326.64"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D23").Value = "2.19"
$ws.Range("E23").Value = "  +3.43%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "10.06"
$ws.Range("E25").Value = "  +17.08%  "
$ws.Range("D26").Value = "65.56"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").Value = "639.48"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.591.45"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0980"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("D30").Value = "0.995"
$ws.Range("E30").Value = "  -17.91%  "
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("E34").Value = "  -2.62%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "1.54"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "152.15"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.368"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "18.67"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").Value = "2.78"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "5.35"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "0.0₆0303"
$ws.Range("E45").Value = "  -26.73%  "
$ws.Range("D46").Value = "153.27"
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "20.38"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").Value = "0.608"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "0.0510"
$ws.Range("E51").Value = "  -1.46%  "
